$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 6819
$ws.Cells.Item(2, 2).Value = 46043.95833333334
$ws.Cells.Item(3, 1).Value = 6756
$ws.Cells.Item(3, 2).Value = 46043.96875
$ws.Cells.Item(4, 1).Value = 6657
$ws.Cells.Item(4, 2).Value = 46043.97916666666
$ws.Cells.Item(5, 1).Value = 6677
$ws.Cells.Item(5, 2).Value = 46043.98958333334
$ws.Cells.Item(6, 1).Value = 6627
$ws.Cells.Item(6, 2).Value = 46044
$ws.Cells.Item(7, 1).Value = 6604
$ws.Cells.Item(7, 2).Value = 46044.01041666666
$ws.Cells.Item(8, 1).Value = 6507
$ws.Cells.Item(8, 2).Value = 46044.02083333334
$ws.Cells.Item(9, 1).Value = 6530
$ws.Cells.Item(9, 2).Value = 46044.03125
$ws.Cells.Item(10, 1).Value = 6440
$ws.Cells.Item(10, 2).Value = 46044.04166666666
$ws.Cells.Item(11, 1).Value = 6493
$ws.Cells.Item(11, 2).Value = 46044.05208333334
$ws.Cells.Item(12, 1).Value = 6426
$ws.Cells.Item(12, 2).Value = 46044.0625
$ws.Cells.Item(13, 1).Value = 6432
$ws.Cells.Item(13, 2).Value = 46044.07291666666
$ws.Cells.Item(14, 1).Value = 6453
$ws.Cells.Item(14, 2).Value = 46044.08333333334
$ws.Cells.Item(15, 1).Value = 6455
$ws.Cells.Item(15, 2).Value = 46044.09375
$ws.Cells.Item(16, 1).Value = 6473
$ws.Cells.Item(16, 2).Value = 46044.10416666666
$ws.Cells.Item(17, 1).Value = 6544
$ws.Cells.Item(17, 2).Value = 46044.11458333334
$ws.Cells.Item(18, 1).Value = 6535
$ws.Cells.Item(18, 2).Value = 46044.125
$ws.Cells.Item(19, 1).Value = 6496
$ws.Cells.Item(19, 2).Value = 46044.13541666666
$ws.Cells.Item(20, 1).Value = 6555
$ws.Cells.Item(20, 2).Value = 46044.14583333334
$ws.Cells.Item(21, 1).Value = 6584
$ws.Cells.Item(21, 2).Value = 46044.15625
$ws.Cells.Item(22, 1).Value = 6646
$ws.Cells.Item(22, 2).Value = 46044.16666666666
$ws.Cells.Item(23, 1).Value = 6557
$ws.Cells.Item(23, 2).Value = 46044.17708333334
$ws.Cells.Item(24, 1).Value = 6601
$ws.Cells.Item(24, 2).Value = 46044.1875
$ws.Cells.Item(25, 1).Value = 6712
$ws.Cells.Item(25, 2).Value = 46044.19791666666
$ws.Cells.Item(26, 1).Value = 6989
$ws.Cells.Item(26, 2).Value = 46044.20833333334
$ws.Cells.Item(27, 1).Value = 7189
$ws.Cells.Item(27, 2).Value = 46044.21875
$ws.Cells.Item(28, 1).Value = 7353
$ws.Cells.Item(28, 2).Value = 46044.22916666666
$ws.Cells.Item(29, 1).Value = 7609
$ws.Cells.Item(29, 2).Value = 46044.23958333334
$ws.Cells.Item(30, 1).Value = 7945
$ws.Cells.Item(30, 2).Value = 46044.25
$ws.Cells.Item(31, 1).Value = 8141
$ws.Cells.Item(31, 2).Value = 46044.26041666666
$ws.Cells.Item(32, 1).Value = 8286
$ws.Cells.Item(32, 2).Value = 46044.27083333334
$ws.Cells.Item(33, 1).Value = 8404
$ws.Cells.Item(33, 2).Value = 46044.28125
$ws.Cells.Item(34, 1).Value = 8562
$ws.Cells.Item(34, 2).Value = 46044.29166666666
$ws.Cells.Item(35, 1).Value = 8703
$ws.Cells.Item(35, 2).Value = 46044.30208333334

$ws.Range("A36:B44").EntireRow.Delete()

$wb.Save()
